# Update attendee counts ("想去人数", column F) on the 展览, 演出 and 全部类型
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 93
$wsExhibit.Range("F3").Value = 11996
$wsExhibit.Range("F4").Value = 34
$wsExhibit.Range("F6").Value = 363
$wsExhibit.Range("F8").Value = 11894
$wsExhibit.Range("F10").Value = 1176
$wsExhibit.Range("F14").Value = 5885
$wsExhibit.Range("F16").Value = 3549

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 575

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 575
$wsAll.Range("F3").Value = 93
$wsAll.Range("F5").Value = 11996
$wsAll.Range("F6").Value = 34
$wsAll.Range("F9").Value = 363
$wsAll.Range("F11").Value = 11894
$wsAll.Range("F13").Value = 1176
$wsAll.Range("F18").Value = 5885
$wsAll.Range("F20").Value = 3549
